$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder rows: swap F:V content between row pairs/cycles sharing the same date ---
# row 24 <= original row 25
$ws.Cells.Item(24,6).Value = "Bilje"
$ws.Cells.Item(24,7).Value = 1
$ws.Cells.Item(24,8).Value = "Tabor Sezana"
$ws.Cells.Item(24,9).Value = 1
$ws.Cells.Item(24,10).Value = 1.52
$ws.Cells.Item(24,11).Value = "20/08/2023 09:00"
$ws.Cells.Item(24,12).Value = 1.47
$ws.Cells.Item(24,13).Value = "20/08/2023 17:19"
$ws.Cells.Item(24,14).Value = 4.26
$ws.Cells.Item(24,15).Value = "20/08/2023 09:00"
$ws.Cells.Item(24,16).Value = 5.59
$ws.Cells.Item(24,17).Value = "20/08/2023 17:19"
$ws.Cells.Item(24,18).Value = 5.03
$ws.Cells.Item(24,19).Value = "20/08/2023 09:00"
$ws.Cells.Item(24,20).Value = 4.39
$ws.Cells.Item(24,21).Value = "20/08/2023 17:25"
$ws.Cells.Item(24,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-tabor-sezana/6ouckJRt/"

# row 25 <= original row 24
$ws.Cells.Item(25,6).Value = "Dravinja"
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = "Ilirija"
$ws.Cells.Item(25,9).Value = 0
$ws.Cells.Item(25,10).Value = 3.12
$ws.Cells.Item(25,11).Value = "20/08/2023 08:59"
$ws.Cells.Item(25,12).Value = 3.08
$ws.Cells.Item(25,13).Value = "20/08/2023 11:20"
$ws.Cells.Item(25,14).Value = 3.27
$ws.Cells.Item(25,15).Value = "20/08/2023 08:59"
$ws.Cells.Item(25,16).Value = 3.4
$ws.Cells.Item(25,17).Value = "20/08/2023 15:35"
$ws.Cells.Item(25,18).Value = 2.15
$ws.Cells.Item(25,19).Value = "20/08/2023 08:59"
$ws.Cells.Item(25,20).Value = 2.12
$ws.Cells.Item(25,21).Value = "20/08/2023 14:26"
$ws.Cells.Item(25,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/dravinja-ilirija/zRWDoyR5/"

# row 29 <= original row 30
$ws.Cells.Item(29,6).Value = "Fuzinar"
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = "Grosuplje"
$ws.Cells.Item(29,9).Value = 4
$ws.Cells.Item(29,10).Value = 3.24
$ws.Cells.Item(29,11).Value = "04/08/2023 05:42"
$ws.Cells.Item(29,12).Value = 4.09
$ws.Cells.Item(29,13).Value = "30/08/2023 15:38"
$ws.Cells.Item(29,14).Value = 3.15
$ws.Cells.Item(29,15).Value = "04/08/2023 05:42"
$ws.Cells.Item(29,16).Value = 3.69
$ws.Cells.Item(29,17).Value = "30/08/2023 15:38"
$ws.Cells.Item(29,18).Value = 2.03
$ws.Cells.Item(29,19).Value = "04/08/2023 05:42"
$ws.Cells.Item(29,20).Value = 1.75
$ws.Cells.Item(29,21).Value = "30/08/2023 15:38"
$ws.Cells.Item(29,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/fuzinar-grosuplje/OMD8pzPE/"

# row 30 <= original row 29
$ws.Cells.Item(30,6).Value = "Ilirija"
$ws.Cells.Item(30,7).Value = 0
$ws.Cells.Item(30,8).Value = "Rudar"
$ws.Cells.Item(30,9).Value = 2
$ws.Cells.Item(30,10).Value = 1.81
$ws.Cells.Item(30,11).Value = "05/08/2023 05:42"
$ws.Cells.Item(30,12).Value = 2.13
$ws.Cells.Item(30,13).Value = "30/08/2023 16:51"
$ws.Cells.Item(30,14).Value = 3.51
$ws.Cells.Item(30,15).Value = "05/08/2023 05:42"
$ws.Cells.Item(30,16).Value = 3.62
$ws.Cells.Item(30,17).Value = "30/08/2023 16:21"
$ws.Cells.Item(30,18).Value = 3.42
$ws.Cells.Item(30,19).Value = "05/08/2023 05:42"
$ws.Cells.Item(30,20).Value = 2.95
$ws.Cells.Item(30,21).Value = "30/08/2023 16:51"
$ws.Cells.Item(30,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-rudar/pv9TMWW7/"

# row 33 <= original row 36
$ws.Cells.Item(33,6).Value = "Primorje"
$ws.Cells.Item(33,7).Value = 3
$ws.Cells.Item(33,8).Value = "Nafta"
$ws.Cells.Item(33,9).Value = 0
$ws.Cells.Item(33,10).Value = 2.32
$ws.Cells.Item(33,11).Value = "01/09/2023 04:43"
$ws.Cells.Item(33,12).Value = 2.33
$ws.Cells.Item(33,13).Value = "02/09/2023 16:27"
$ws.Cells.Item(33,14).Value = 3.25
$ws.Cells.Item(33,15).Value = "01/09/2023 04:43"
$ws.Cells.Item(33,16).Value = 3.38
$ws.Cells.Item(33,17).Value = "02/09/2023 16:27"
$ws.Cells.Item(33,18).Value = 2.63
$ws.Cells.Item(33,19).Value = "01/09/2023 04:43"
$ws.Cells.Item(33,20).Value = 2.77
$ws.Cells.Item(33,21).Value = "02/09/2023 16:27"
$ws.Cells.Item(33,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/primorje-nafta/GKUtkAXc/"

# row 34 <= original row 35
$ws.Cells.Item(34,6).Value = "Tabor Sezana"
$ws.Cells.Item(34,7).Value = 1
$ws.Cells.Item(34,8).Value = "Rudar"
$ws.Cells.Item(34,9).Value = 1
$ws.Cells.Item(34,10).Value = 4.32
$ws.Cells.Item(34,11).Value = "02/09/2023 14:10"
$ws.Cells.Item(34,12).Value = 3.42
$ws.Cells.Item(34,13).Value = "02/09/2023 16:28"
$ws.Cells.Item(34,14).Value = 3.88
$ws.Cells.Item(34,15).Value = "02/09/2023 14:10"
$ws.Cells.Item(34,16).Value = 3.99
$ws.Cells.Item(34,17).Value = "02/09/2023 16:28"
$ws.Cells.Item(34,18).Value = 1.66
$ws.Cells.Item(34,19).Value = "02/09/2023 14:10"
$ws.Cells.Item(34,20).Value = 1.85
$ws.Cells.Item(34,21).Value = "02/09/2023 16:28"
$ws.Cells.Item(34,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tabor-sezana-rudar/UNYplUn4/"

# row 35 <= original row 34
$ws.Cells.Item(35,6).Value = "Fuzinar"
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(35,8).Value = "Jadran Dekani"
$ws.Cells.Item(35,9).Value = 3
$ws.Cells.Item(35,10).Value = 2.35
$ws.Cells.Item(35,11).Value = "01/09/2023 04:43"
$ws.Cells.Item(35,12).Value = 2.84
$ws.Cells.Item(35,13).Value = "02/09/2023 16:21"
$ws.Cells.Item(35,14).Value = 3.12
$ws.Cells.Item(35,15).Value = "01/09/2023 04:43"
$ws.Cells.Item(35,16).Value = 3.36
$ws.Cells.Item(35,17).Value = "02/09/2023 16:21"
$ws.Cells.Item(35,18).Value = 2.68
$ws.Cells.Item(35,19).Value = "01/09/2023 04:43"
$ws.Cells.Item(35,20).Value = 2.3
$ws.Cells.Item(35,21).Value = "02/09/2023 16:21"
$ws.Cells.Item(35,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/fuzinar-jadran-dekani/YeuxjjIi/"

# row 36 <= original row 33
$ws.Cells.Item(36,6).Value = "ND Gorica"
$ws.Cells.Item(36,7).Value = 3
$ws.Cells.Item(36,8).Value = "Dravinja"
$ws.Cells.Item(36,9).Value = 0
$ws.Cells.Item(36,10).Value = 1.58
$ws.Cells.Item(36,11).Value = "02/09/2023 14:10"
$ws.Cells.Item(36,12).Value = 1.44
$ws.Cells.Item(36,13).Value = "02/09/2023 14:41"
$ws.Cells.Item(36,14).Value = 3.99
$ws.Cells.Item(36,15).Value = "02/09/2023 14:10"
$ws.Cells.Item(36,16).Value = 4.38
$ws.Cells.Item(36,17).Value = "02/09/2023 14:41"
$ws.Cells.Item(36,18).Value = 4.75
$ws.Cells.Item(36,19).Value = "02/09/2023 14:10"
$ws.Cells.Item(36,20).Value = 6.18
$ws.Cells.Item(36,21).Value = "02/09/2023 14:41"
$ws.Cells.Item(36,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-dravinja/f5tYjW2o/"

# row 44 <= original row 46
$ws.Cells.Item(44,6).Value = "Nafta"
$ws.Cells.Item(44,7).Value = 0
$ws.Cells.Item(44,8).Value = "Tabor Sezana"
$ws.Cells.Item(44,9).Value = 0
$ws.Cells.Item(44,10).Value = 1.19
$ws.Cells.Item(44,11).Value = "09/09/2023 13:42"
$ws.Cells.Item(44,12).Value = 1.29
$ws.Cells.Item(44,13).Value = "09/09/2023 16:23"
$ws.Cells.Item(44,14).Value = 6.85
$ws.Cells.Item(44,15).Value = "09/09/2023 13:42"
$ws.Cells.Item(44,16).Value = 5.9
$ws.Cells.Item(44,17).Value = "09/09/2023 16:23"
$ws.Cells.Item(44,18).Value = 8.869999999999999
$ws.Cells.Item(44,19).Value = "09/09/2023 13:42"
$ws.Cells.Item(44,20).Value = 7.25
$ws.Cells.Item(44,21).Value = "09/09/2023 16:23"
$ws.Cells.Item(44,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nafta-tabor-sezana/hMwQvQ9j/"

# row 46 <= original row 47
$ws.Cells.Item(46,6).Value = "Rudar"
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(46,8).Value = "Beltinci"
$ws.Cells.Item(46,9).Value = 1
$ws.Cells.Item(46,10).Value = 2.41
$ws.Cells.Item(46,11).Value = "08/09/2023 04:42"
$ws.Cells.Item(46,12).Value = 3.4
$ws.Cells.Item(46,13).Value = "09/09/2023 16:15"
$ws.Cells.Item(46,14).Value = 3.23
$ws.Cells.Item(46,15).Value = "08/09/2023 04:42"
$ws.Cells.Item(46,16).Value = 3.56
$ws.Cells.Item(46,17).Value = "09/09/2023 16:14"
$ws.Cells.Item(46,18).Value = 2.54
$ws.Cells.Item(46,19).Value = "08/09/2023 04:42"
$ws.Cells.Item(46,20).Value = 1.91
$ws.Cells.Item(46,21).Value = "09/09/2023 16:15"
$ws.Cells.Item(46,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-beltinci/vJZLu6fp/"

# row 47 <= original row 44
$ws.Cells.Item(47,6).Value = "Bilje"
$ws.Cells.Item(47,7).Value = 2
$ws.Cells.Item(47,8).Value = "Tolmin"
$ws.Cells.Item(47,9).Value = 1
$ws.Cells.Item(47,10).Value = 1.54
$ws.Cells.Item(47,11).Value = "08/09/2023 04:42"
$ws.Cells.Item(47,12).Value = 1.49
$ws.Cells.Item(47,13).Value = "09/09/2023 16:22"
$ws.Cells.Item(47,14).Value = 3.85
$ws.Cells.Item(47,15).Value = "08/09/2023 04:42"
$ws.Cells.Item(47,16).Value = 4.38
$ws.Cells.Item(47,17).Value = "09/09/2023 16:22"
$ws.Cells.Item(47,18).Value = 4.49
$ws.Cells.Item(47,19).Value = "08/09/2023 04:42"
$ws.Cells.Item(47,20).Value = 5.37
$ws.Cells.Item(47,21).Value = "09/09/2023 16:22"
$ws.Cells.Item(47,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-tolmin/MwmfpnnT/"

# row 66 <= original row 68
$ws.Cells.Item(66,6).Value = "Tabor Sezana"
$ws.Cells.Item(66,7).Value = 0
$ws.Cells.Item(66,8).Value = "NK Krka"
$ws.Cells.Item(66,9).Value = 2
$ws.Cells.Item(66,10).Value = 3.75
$ws.Cells.Item(66,11).Value = "22/09/2023 03:13"
$ws.Cells.Item(66,12).Value = 4.99
$ws.Cells.Item(66,13).Value = "23/09/2023 15:47"
$ws.Cells.Item(66,14).Value = 3.56
$ws.Cells.Item(66,15).Value = "22/09/2023 03:13"
$ws.Cells.Item(66,16).Value = 4.23
$ws.Cells.Item(66,17).Value = "23/09/2023 15:47"
$ws.Cells.Item(66,18).Value = 1.72
$ws.Cells.Item(66,19).Value = "22/09/2023 03:13"
$ws.Cells.Item(66,20).Value = 1.54
$ws.Cells.Item(66,21).Value = "23/09/2023 15:47"
$ws.Cells.Item(66,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tabor-sezana-nk-krka/pf1lPsFR/"

# row 68 <= original row 66
$ws.Cells.Item(68,6).Value = "Fuzinar"
$ws.Cells.Item(68,7).Value = 3
$ws.Cells.Item(68,8).Value = "Ilirija"
$ws.Cells.Item(68,9).Value = 2
$ws.Cells.Item(68,10).Value = 2.44
$ws.Cells.Item(68,11).Value = "22/09/2023 03:13"
$ws.Cells.Item(68,12).Value = 2.45
$ws.Cells.Item(68,13).Value = "23/09/2023 15:58"
$ws.Cells.Item(68,14).Value = 3.32
$ws.Cells.Item(68,15).Value = "22/09/2023 03:13"
$ws.Cells.Item(68,16).Value = 3.7
$ws.Cells.Item(68,17).Value = "23/09/2023 15:59"
$ws.Cells.Item(68,18).Value = 2.45
$ws.Cells.Item(68,19).Value = "22/09/2023 03:13"
$ws.Cells.Item(68,20).Value = 2.46
$ws.Cells.Item(68,21).Value = "23/09/2023 15:58"
$ws.Cells.Item(68,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/fuzinar-ilirija/dY8ySqV8/"

# row 69 <= original row 70
$ws.Cells.Item(69,6).Value = "Tolmin"
$ws.Cells.Item(69,7).Value = 2
$ws.Cells.Item(69,8).Value = "Jadran Dekani"
$ws.Cells.Item(69,9).Value = 2
$ws.Cells.Item(69,10).Value = 3.15
$ws.Cells.Item(69,11).Value = "23/09/2023 03:12"
$ws.Cells.Item(69,12).Value = 3.52
$ws.Cells.Item(69,13).Value = "24/09/2023 15:42"
$ws.Cells.Item(69,14).Value = 3.16
$ws.Cells.Item(69,15).Value = "23/09/2023 03:12"
$ws.Cells.Item(69,16).Value = 3.17
$ws.Cells.Item(69,17).Value = "24/09/2023 15:41"
$ws.Cells.Item(69,18).Value = 2.06
$ws.Cells.Item(69,19).Value = "23/09/2023 03:12"
$ws.Cells.Item(69,20).Value = 2.06
$ws.Cells.Item(69,21).Value = "24/09/2023 15:41"
$ws.Cells.Item(69,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tolmin-jadran-dekani/GxhKJLxk/"

# row 70 <= original row 71
$ws.Cells.Item(70,6).Value = "Triglav"
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = "ND Gorica"
$ws.Cells.Item(70,9).Value = 3
$ws.Cells.Item(70,10).Value = 3.01
$ws.Cells.Item(70,11).Value = "23/09/2023 03:12"
$ws.Cells.Item(70,12).Value = 3.77
$ws.Cells.Item(70,13).Value = "24/09/2023 15:32"
$ws.Cells.Item(70,14).Value = 3.29
$ws.Cells.Item(70,15).Value = "23/09/2023 03:12"
$ws.Cells.Item(70,16).Value = 3.27
$ws.Cells.Item(70,17).Value = "24/09/2023 15:39"
$ws.Cells.Item(70,18).Value = 2.03
$ws.Cells.Item(70,19).Value = "23/09/2023 03:12"
$ws.Cells.Item(70,20).Value = 1.94
$ws.Cells.Item(70,21).Value = "24/09/2023 15:39"
$ws.Cells.Item(70,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/triglav-nd-gorica/IsAXSPF2/"

# row 71 <= original row 69
$ws.Cells.Item(71,6).Value = "Beltinci"
$ws.Cells.Item(71,7).Value = 3
$ws.Cells.Item(71,8).Value = "Dravinja"
$ws.Cells.Item(71,9).Value = 0
$ws.Cells.Item(71,10).Value = 1.4
$ws.Cells.Item(71,11).Value = "23/09/2023 03:12"
$ws.Cells.Item(71,12).Value = 1.38
$ws.Cells.Item(71,13).Value = "24/09/2023 15:43"
$ws.Cells.Item(71,14).Value = 4.36
$ws.Cells.Item(71,15).Value = "23/09/2023 03:12"
$ws.Cells.Item(71,16).Value = 4.76
$ws.Cells.Item(71,17).Value = "24/09/2023 15:54"
$ws.Cells.Item(71,18).Value = 5.39
$ws.Cells.Item(71,19).Value = "23/09/2023 03:12"
$ws.Cells.Item(71,20).Value = 6.65
$ws.Cells.Item(71,21).Value = "24/09/2023 15:54"
$ws.Cells.Item(71,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/beltinci-dravinja/n3mGK1Nr/"

# row 73 <= original row 74
$ws.Cells.Item(73,6).Value = "Ilirija"
$ws.Cells.Item(73,7).Value = 0
$ws.Cells.Item(73,8).Value = "Primorje"
$ws.Cells.Item(73,9).Value = 0
$ws.Cells.Item(73,10).Value = 3.72
$ws.Cells.Item(73,11).Value = "28/09/2023 02:42"
$ws.Cells.Item(73,12).Value = 4.94
$ws.Cells.Item(73,13).Value = "29/09/2023 15:20"
$ws.Cells.Item(73,14).Value = 3.4
$ws.Cells.Item(73,15).Value = "28/09/2023 02:42"
$ws.Cells.Item(73,16).Value = 3.8
$ws.Cells.Item(73,17).Value = "29/09/2023 15:20"
$ws.Cells.Item(73,18).Value = 1.79
$ws.Cells.Item(73,19).Value = "28/09/2023 02:42"
$ws.Cells.Item(73,20).Value = 1.61
$ws.Cells.Item(73,21).Value = "29/09/2023 15:20"
$ws.Cells.Item(73,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-primorje/f1UjKhTa/"

# row 74 <= original row 73
$ws.Cells.Item(74,6).Value = "ND Gorica"
$ws.Cells.Item(74,7).Value = 2
$ws.Cells.Item(74,8).Value = "Fuzinar"
$ws.Cells.Item(74,9).Value = 1
$ws.Cells.Item(74,10).Value = 1.29
$ws.Cells.Item(74,11).Value = "28/09/2023 02:42"
$ws.Cells.Item(74,12).Value = 1.25
$ws.Cells.Item(74,13).Value = "29/09/2023 13:45"
$ws.Cells.Item(74,14).Value = 4.98
$ws.Cells.Item(74,15).Value = "28/09/2023 02:42"
$ws.Cells.Item(74,16).Value = 6.11
$ws.Cells.Item(74,17).Value = "29/09/2023 15:29"
$ws.Cells.Item(74,18).Value = 6.88
$ws.Cells.Item(74,19).Value = "28/09/2023 02:42"
$ws.Cells.Item(74,20).Value = 7.81
$ws.Cells.Item(74,21).Value = "29/09/2023 15:29"
$ws.Cells.Item(74,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-fuzinar/xQJeJCr6/"

# row 76 <= original row 77
$ws.Cells.Item(76,6).Value = "NK Bistrica"
$ws.Cells.Item(76,7).Value = 1
$ws.Cells.Item(76,8).Value = "Tabor Sezana"
$ws.Cells.Item(76,9).Value = 0
$ws.Cells.Item(76,10).Value = 1.41
$ws.Cells.Item(76,11).Value = "29/09/2023 02:42"
$ws.Cells.Item(76,12).Value = 1.34
$ws.Cells.Item(76,13).Value = "30/09/2023 15:24"
$ws.Cells.Item(76,14).Value = 4.29
$ws.Cells.Item(76,15).Value = "29/09/2023 02:42"
$ws.Cells.Item(76,16).Value = 5.2
$ws.Cells.Item(76,17).Value = "30/09/2023 15:29"
$ws.Cells.Item(76,18).Value = 5.33
$ws.Cells.Item(76,19).Value = "29/09/2023 02:42"
$ws.Cells.Item(76,20).Value = 6.9
$ws.Cells.Item(76,21).Value = "30/09/2023 15:29"
$ws.Cells.Item(76,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bistrica-tabor-sezana/tbQnLYDg/"

# row 77 <= original row 78
$ws.Cells.Item(77,6).Value = "Jadran Dekani"
$ws.Cells.Item(77,7).Value = 1
$ws.Cells.Item(77,8).Value = "Grosuplje"
$ws.Cells.Item(77,9).Value = 1
$ws.Cells.Item(77,10).Value = 3.27
$ws.Cells.Item(77,11).Value = "29/09/2023 02:42"
$ws.Cells.Item(77,12).Value = 4
$ws.Cells.Item(77,13).Value = "30/09/2023 15:20"
$ws.Cells.Item(77,14).Value = 3.13
$ws.Cells.Item(77,15).Value = "29/09/2023 02:42"
$ws.Cells.Item(77,16).Value = 3.26
$ws.Cells.Item(77,17).Value = "30/09/2023 15:20"
$ws.Cells.Item(77,18).Value = 1.99
$ws.Cells.Item(77,19).Value = "29/09/2023 02:42"
$ws.Cells.Item(77,20).Value = 1.89
$ws.Cells.Item(77,21).Value = "30/09/2023 15:20"
$ws.Cells.Item(77,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/jadran-dekani-grosuplje/rkXIvNEE/"

# row 78 <= original row 79
$ws.Cells.Item(78,6).Value = "Nafta"
$ws.Cells.Item(78,7).Value = 3
$ws.Cells.Item(78,8).Value = "Bilje"
$ws.Cells.Item(78,9).Value = 1
$ws.Cells.Item(78,10).Value = 1.72
$ws.Cells.Item(78,11).Value = "29/09/2023 02:42"
$ws.Cells.Item(78,12).Value = 1.78
$ws.Cells.Item(78,13).Value = "30/09/2023 15:23"
$ws.Cells.Item(78,14).Value = 3.77
$ws.Cells.Item(78,15).Value = "29/09/2023 02:42"
$ws.Cells.Item(78,16).Value = 4.06
$ws.Cells.Item(78,17).Value = "30/09/2023 15:23"
$ws.Cells.Item(78,18).Value = 3.65
$ws.Cells.Item(78,19).Value = "29/09/2023 02:42"
$ws.Cells.Item(78,20).Value = 3.61
$ws.Cells.Item(78,21).Value = "30/09/2023 15:23"
$ws.Cells.Item(78,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nafta-bilje/4bYEu3a8/"

# row 79 <= original row 76
$ws.Cells.Item(79,6).Value = "Rudar"
$ws.Cells.Item(79,7).Value = 0
$ws.Cells.Item(79,8).Value = "Triglav"
$ws.Cells.Item(79,9).Value = 2
$ws.Cells.Item(79,10).Value = 2.29
$ws.Cells.Item(79,11).Value = "29/09/2023 02:42"
$ws.Cells.Item(79,12).Value = 2.25
$ws.Cells.Item(79,13).Value = "30/09/2023 15:29"
$ws.Cells.Item(79,14).Value = 3.2
$ws.Cells.Item(79,15).Value = "29/09/2023 02:42"
$ws.Cells.Item(79,16).Value = 3.37
$ws.Cells.Item(79,17).Value = "30/09/2023 15:25"
$ws.Cells.Item(79,18).Value = 2.64
$ws.Cells.Item(79,19).Value = "29/09/2023 02:42"
$ws.Cells.Item(79,20).Value = 2.9
$ws.Cells.Item(79,21).Value = "30/09/2023 15:24"
$ws.Cells.Item(79,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-triglav/jepXGJM7/"

# row 82 <= original row 83
$ws.Cells.Item(82,6).Value = "Beltinci"
$ws.Cells.Item(82,7).Value = 3
$ws.Cells.Item(82,8).Value = "NK Bistrica"
$ws.Cells.Item(82,9).Value = 0
$ws.Cells.Item(82,10).Value = 1.71
$ws.Cells.Item(82,11).Value = "06/10/2023 02:42"
$ws.Cells.Item(82,12).Value = 1.52
$ws.Cells.Item(82,13).Value = "07/10/2023 10:57"
$ws.Cells.Item(82,14).Value = 3.56
$ws.Cells.Item(82,15).Value = "06/10/2023 02:42"
$ws.Cells.Item(82,16).Value = 4.11
$ws.Cells.Item(82,17).Value = "07/10/2023 13:35"
$ws.Cells.Item(82,18).Value = 3.79
$ws.Cells.Item(82,19).Value = "06/10/2023 02:42"
$ws.Cells.Item(82,20).Value = 5.27
$ws.Cells.Item(82,21).Value = "07/10/2023 10:57"
$ws.Cells.Item(82,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/beltinci-bistrica/Qo0juf5P/"

# row 83 <= original row 84
$ws.Cells.Item(83,6).Value = "Bilje"
$ws.Cells.Item(83,7).Value = 2
$ws.Cells.Item(83,8).Value = "Jadran Dekani"
$ws.Cells.Item(83,9).Value = 2
$ws.Cells.Item(83,10).Value = 2.04
$ws.Cells.Item(83,11).Value = "06/10/2023 02:42"
$ws.Cells.Item(83,12).Value = 2.34
$ws.Cells.Item(83,13).Value = "07/10/2023 15:15"
$ws.Cells.Item(83,14).Value = 3.22
$ws.Cells.Item(83,15).Value = "06/10/2023 02:42"
$ws.Cells.Item(83,16).Value = 3.34
$ws.Cells.Item(83,17).Value = "07/10/2023 15:15"
$ws.Cells.Item(83,18).Value = 3.06
$ws.Cells.Item(83,19).Value = "06/10/2023 02:42"
$ws.Cells.Item(83,20).Value = 2.8
$ws.Cells.Item(83,21).Value = "07/10/2023 15:15"
$ws.Cells.Item(83,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-jadran-dekani/YikQYirg/"

# row 84 <= original row 82
$ws.Cells.Item(84,6).Value = "Grosuplje"
$ws.Cells.Item(84,7).Value = 1
$ws.Cells.Item(84,8).Value = "Dravinja"
$ws.Cells.Item(84,9).Value = 2
$ws.Cells.Item(84,10).Value = 1.34
$ws.Cells.Item(84,11).Value = "06/10/2023 02:42"
$ws.Cells.Item(84,12).Value = 1.29
$ws.Cells.Item(84,13).Value = "07/10/2023 15:03"
$ws.Cells.Item(84,14).Value = 4.49
$ws.Cells.Item(84,15).Value = "06/10/2023 02:42"
$ws.Cells.Item(84,16).Value = 5.25
$ws.Cells.Item(84,17).Value = "07/10/2023 15:19"
$ws.Cells.Item(84,18).Value = 6.24
$ws.Cells.Item(84,19).Value = "06/10/2023 02:42"
$ws.Cells.Item(84,20).Value = 8.34
$ws.Cells.Item(84,21).Value = "07/10/2023 15:19"
$ws.Cells.Item(84,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-dravinja/f1jMZXSn/"

# row 85 <= original row 87
$ws.Cells.Item(85,6).Value = "Tolmin"
$ws.Cells.Item(85,7).Value = 1
$ws.Cells.Item(85,8).Value = "NK Krka"
$ws.Cells.Item(85,9).Value = 3
$ws.Cells.Item(85,10).Value = 3.75
$ws.Cells.Item(85,11).Value = "07/10/2023 02:42"
$ws.Cells.Item(85,12).Value = 2.86
$ws.Cells.Item(85,13).Value = "08/10/2023 15:27"
$ws.Cells.Item(85,14).Value = 3.48
$ws.Cells.Item(85,15).Value = "07/10/2023 02:42"
$ws.Cells.Item(85,16).Value = 3.82
$ws.Cells.Item(85,17).Value = "08/10/2023 15:28"
$ws.Cells.Item(85,18).Value = 1.74
$ws.Cells.Item(85,19).Value = "07/10/2023 02:42"
$ws.Cells.Item(85,20).Value = 2.11
$ws.Cells.Item(85,21).Value = "08/10/2023 15:27"
$ws.Cells.Item(85,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tolmin-nk-krka/0viIzEDt/"

# row 86 <= original row 85
$ws.Cells.Item(86,6).Value = "Rudar"
$ws.Cells.Item(86,7).Value = 1
$ws.Cells.Item(86,8).Value = "Nafta"
$ws.Cells.Item(86,9).Value = 5
$ws.Cells.Item(86,10).Value = 2.62
$ws.Cells.Item(86,11).Value = "07/10/2023 02:42"
$ws.Cells.Item(86,12).Value = 3.04
$ws.Cells.Item(86,13).Value = "08/10/2023 15:01"
$ws.Cells.Item(86,14).Value = 3.22
$ws.Cells.Item(86,15).Value = "07/10/2023 02:42"
$ws.Cells.Item(86,16).Value = 3.36
$ws.Cells.Item(86,17).Value = "08/10/2023 15:01"
$ws.Cells.Item(86,18).Value = 2.29
$ws.Cells.Item(86,19).Value = "07/10/2023 02:42"
$ws.Cells.Item(86,20).Value = 2.18
$ws.Cells.Item(86,21).Value = "08/10/2023 15:01"
$ws.Cells.Item(86,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-nafta/jTqVXBca/"

# row 87 <= original row 86
$ws.Cells.Item(87,6).Value = "Tabor Sezana"
$ws.Cells.Item(87,7).Value = 0
$ws.Cells.Item(87,8).Value = "Ilirija"
$ws.Cells.Item(87,9).Value = 3
$ws.Cells.Item(87,10).Value = 2.97
$ws.Cells.Item(87,11).Value = "07/10/2023 02:42"
$ws.Cells.Item(87,12).Value = 3.94
$ws.Cells.Item(87,13).Value = "08/10/2023 13:47"
$ws.Cells.Item(87,14).Value = 3.28
$ws.Cells.Item(87,15).Value = "07/10/2023 02:42"
$ws.Cells.Item(87,16).Value = 3.92
$ws.Cells.Item(87,17).Value = "08/10/2023 15:01"
$ws.Cells.Item(87,18).Value = 2.05
$ws.Cells.Item(87,19).Value = "07/10/2023 02:42"
$ws.Cells.Item(87,20).Value = 1.72
$ws.Cells.Item(87,21).Value = "08/10/2023 13:47"
$ws.Cells.Item(87,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tabor-sezana-ilirija/OtM7GASO/"

# row 98 <= original row 99
$ws.Cells.Item(98,6).Value = "Grosuplje"
$ws.Cells.Item(98,7).Value = 1
$ws.Cells.Item(98,8).Value = "NK Bistrica"
$ws.Cells.Item(98,9).Value = 1
$ws.Cells.Item(98,10).Value = 1.75
$ws.Cells.Item(98,11).Value = "19/10/2023 02:12"
$ws.Cells.Item(98,12).Value = 1.75
$ws.Cells.Item(98,13).Value = "20/10/2023 14:57"
$ws.Cells.Item(98,14).Value = 3.44
$ws.Cells.Item(98,15).Value = "19/10/2023 02:12"
$ws.Cells.Item(98,16).Value = 3.56
$ws.Cells.Item(98,17).Value = "20/10/2023 14:57"
$ws.Cells.Item(98,18).Value = 3.74
$ws.Cells.Item(98,19).Value = "19/10/2023 02:12"
$ws.Cells.Item(98,20).Value = 4.32
$ws.Cells.Item(98,21).Value = "20/10/2023 14:57"
$ws.Cells.Item(98,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-bistrica/hMePsmRN/"

# row 99 <= original row 98
$ws.Cells.Item(99,6).Value = "Bilje"
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = "NK Krka"
$ws.Cells.Item(99,9).Value = 2
$ws.Cells.Item(99,10).Value = 2.55
$ws.Cells.Item(99,11).Value = "19/10/2023 02:12"
$ws.Cells.Item(99,12).Value = 2.55
$ws.Cells.Item(99,13).Value = "20/10/2023 14:58"
$ws.Cells.Item(99,14).Value = 3.22
$ws.Cells.Item(99,15).Value = "19/10/2023 02:12"
$ws.Cells.Item(99,16).Value = 3.76
$ws.Cells.Item(99,17).Value = "20/10/2023 14:58"
$ws.Cells.Item(99,18).Value = 2.36
$ws.Cells.Item(99,19).Value = "19/10/2023 02:12"
$ws.Cells.Item(99,20).Value = 2.35
$ws.Cells.Item(99,21).Value = "20/10/2023 14:58"
$ws.Cells.Item(99,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-nk-krka/nqeTt7tU/"

# row 106 <= original row 107
$ws.Cells.Item(106,6).Value = "Bilje"
$ws.Cells.Item(106,7).Value = 5
$ws.Cells.Item(106,8).Value = "Fuzinar"
$ws.Cells.Item(106,9).Value = 1
$ws.Cells.Item(106,10).Value = 1.89
$ws.Cells.Item(106,11).Value = "23/10/2023 02:12"
$ws.Cells.Item(106,12).Value = 2
$ws.Cells.Item(106,13).Value = "24/10/2023 14:56"
$ws.Cells.Item(106,14).Value = 3.58
$ws.Cells.Item(106,15).Value = "23/10/2023 02:12"
$ws.Cells.Item(106,16).Value = 3.68
$ws.Cells.Item(106,17).Value = "24/10/2023 14:58"
$ws.Cells.Item(106,18).Value = 3.12
$ws.Cells.Item(106,19).Value = "23/10/2023 02:12"
$ws.Cells.Item(106,20).Value = 3.19
$ws.Cells.Item(106,21).Value = "24/10/2023 14:58"
$ws.Cells.Item(106,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-fuzinar/OpdeyWo1/"

# row 107 <= original row 106
$ws.Cells.Item(107,6).Value = "Grosuplje"
$ws.Cells.Item(107,7).Value = 1
$ws.Cells.Item(107,8).Value = "Primorje"
$ws.Cells.Item(107,9).Value = 0
$ws.Cells.Item(107,10).Value = 2.62
$ws.Cells.Item(107,11).Value = "08/08/2023 04:42"
$ws.Cells.Item(107,12).Value = 2.54
$ws.Cells.Item(107,13).Value = "24/10/2023 14:54"
$ws.Cells.Item(107,14).Value = 3.06
$ws.Cells.Item(107,15).Value = "08/08/2023 04:42"
$ws.Cells.Item(107,16).Value = 3.05
$ws.Cells.Item(107,17).Value = "24/10/2023 14:59"
$ws.Cells.Item(107,18).Value = 2.44
$ws.Cells.Item(107,19).Value = "08/08/2023 04:42"
$ws.Cells.Item(107,20).Value = 2.75
$ws.Cells.Item(107,21).Value = "24/10/2023 14:54"
$ws.Cells.Item(107,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-primorje/pfcixCWf/"

# row 117 <= original row 120
$ws.Cells.Item(117,6).Value = "NK Krka"
$ws.Cells.Item(117,7).Value = 4
$ws.Cells.Item(117,8).Value = "Rudar"
$ws.Cells.Item(117,9).Value = 5
$ws.Cells.Item(117,10).Value = 1.9
$ws.Cells.Item(117,11).Value = "28/10/2023 03:12"
$ws.Cells.Item(117,12).Value = 2.06
$ws.Cells.Item(117,13).Value = "29/10/2023 13:33"
$ws.Cells.Item(117,14).Value = 3.4
$ws.Cells.Item(117,15).Value = "28/10/2023 03:12"
$ws.Cells.Item(117,16).Value = 3.61
$ws.Cells.Item(117,17).Value = "29/10/2023 13:33"
$ws.Cells.Item(117,18).Value = 3.25
$ws.Cells.Item(117,19).Value = "28/10/2023 03:12"
$ws.Cells.Item(117,20).Value = 3.09
$ws.Cells.Item(117,21).Value = "29/10/2023 13:33"
$ws.Cells.Item(117,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nk-krka-rudar/d8skb8B4/"

# row 120 <= original row 117
$ws.Cells.Item(120,6).Value = "Dravinja"
$ws.Cells.Item(120,7).Value = 0
$ws.Cells.Item(120,8).Value = "Nafta"
$ws.Cells.Item(120,9).Value = 1
$ws.Cells.Item(120,10).Value = 3.82
$ws.Cells.Item(120,11).Value = "28/10/2023 03:12"
$ws.Cells.Item(120,12).Value = 3.46
$ws.Cells.Item(120,13).Value = "29/10/2023 13:08"
$ws.Cells.Item(120,14).Value = 3.57
$ws.Cells.Item(120,15).Value = "28/10/2023 03:12"
$ws.Cells.Item(120,16).Value = 3.62
$ws.Cells.Item(120,17).Value = "29/10/2023 13:08"
$ws.Cells.Item(120,18).Value = 1.7
$ws.Cells.Item(120,19).Value = "28/10/2023 03:12"
$ws.Cells.Item(120,20).Value = 1.92
$ws.Cells.Item(120,21).Value = "29/10/2023 13:08"
$ws.Cells.Item(120,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/dravinja-nafta/IHroaldb/"

# row 127 <= original row 128
$ws.Cells.Item(127,6).Value = "NK Bistrica"
$ws.Cells.Item(127,7).Value = 5
$ws.Cells.Item(127,8).Value = "Rudar"
$ws.Cells.Item(127,9).Value = 1
$ws.Cells.Item(127,10).Value = 1.95
$ws.Cells.Item(127,11).Value = "05/11/2023 02:12"
$ws.Cells.Item(127,12).Value = 2.17
$ws.Cells.Item(127,13).Value = "06/11/2023 13:56"
$ws.Cells.Item(127,14).Value = 3.4
$ws.Cells.Item(127,15).Value = "05/11/2023 02:12"
$ws.Cells.Item(127,16).Value = 3.44
$ws.Cells.Item(127,17).Value = "06/11/2023 13:56"
$ws.Cells.Item(127,18).Value = 3.1
$ws.Cells.Item(127,19).Value = "05/11/2023 02:12"
$ws.Cells.Item(127,20).Value = 3
$ws.Cells.Item(127,21).Value = "06/11/2023 13:56"
$ws.Cells.Item(127,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bistrica-rudar/Q3NT95AG/"

# row 128 <= original row 127
$ws.Cells.Item(128,6).Value = "Ilirija"
$ws.Cells.Item(128,7).Value = 0
$ws.Cells.Item(128,8).Value = "Bilje"
$ws.Cells.Item(128,9).Value = 3
$ws.Cells.Item(128,10).Value = 2.03
$ws.Cells.Item(128,11).Value = "05/11/2023 02:12"
$ws.Cells.Item(128,12).Value = 2.04
$ws.Cells.Item(128,13).Value = "06/11/2023 13:57"
$ws.Cells.Item(128,14).Value = 3.41
$ws.Cells.Item(128,15).Value = "05/11/2023 02:12"
$ws.Cells.Item(128,16).Value = 3.66
$ws.Cells.Item(128,17).Value = "06/11/2023 13:57"
$ws.Cells.Item(128,18).Value = 2.92
$ws.Cells.Item(128,19).Value = "05/11/2023 02:12"
$ws.Cells.Item(128,20).Value = 3.12
$ws.Cells.Item(128,21).Value = "06/11/2023 13:57"
$ws.Cells.Item(128,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-bilje/dUIPAoeA/"

# --- Append 4 new match rows (130-133) ---

# new row 130
$ws.Cells.Item(129,1).Copy()
$ws.Cells.Item(130,1).PasteSpecial(-4122)
$ws.Cells.Item(129,5).Copy()
$ws.Cells.Item(130,5).PasteSpecial(-4122)
$ws.Cells.Item(130,1).Value = 129
$ws.Cells.Item(130,2).Value = "slovenia"
$ws.Cells.Item(130,3).Value = "2-snl"
$ws.Cells.Item(130,4).Value = "2023-2024"
$ws.Cells.Item(130,5).Value = 45241.58333333334
$ws.Cells.Item(130,6).Value = "Dravinja"
$ws.Cells.Item(130,7).Value = 2
$ws.Cells.Item(130,8).Value = "Triglav"
$ws.Cells.Item(130,9).Value = 1
$ws.Cells.Item(130,10).Value = 2.58
$ws.Cells.Item(130,11).Value = "10/11/2023 02:13"
$ws.Cells.Item(130,12).Value = 2.8
$ws.Cells.Item(130,13).Value = "11/11/2023 13:39"
$ws.Cells.Item(130,14).Value = 3.09
$ws.Cells.Item(130,15).Value = "10/11/2023 02:13"
$ws.Cells.Item(130,16).Value = 3.34
$ws.Cells.Item(130,17).Value = "11/11/2023 13:39"
$ws.Cells.Item(130,18).Value = 2.4
$ws.Cells.Item(130,19).Value = "10/11/2023 02:13"
$ws.Cells.Item(130,20).Value = 2.34
$ws.Cells.Item(130,21).Value = "11/11/2023 13:39"
$ws.Cells.Item(130,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/dravinja-triglav/dSw6jH9e/"

# new row 131
$ws.Cells.Item(129,1).Copy()
$ws.Cells.Item(131,1).PasteSpecial(-4122)
$ws.Cells.Item(129,5).Copy()
$ws.Cells.Item(131,5).PasteSpecial(-4122)
$ws.Cells.Item(131,1).Value = 130
$ws.Cells.Item(131,2).Value = "slovenia"
$ws.Cells.Item(131,3).Value = "2-snl"
$ws.Cells.Item(131,4).Value = "2023-2024"
$ws.Cells.Item(131,5).Value = 45241.58333333334
$ws.Cells.Item(131,6).Value = "Grosuplje"
$ws.Cells.Item(131,7).Value = 2
$ws.Cells.Item(131,8).Value = "Fuzinar"
$ws.Cells.Item(131,9).Value = 0
$ws.Cells.Item(131,10).Value = 1.39
$ws.Cells.Item(131,11).Value = "10/11/2023 02:13"
$ws.Cells.Item(131,12).Value = 1.56
$ws.Cells.Item(131,13).Value = "11/11/2023 13:56"
$ws.Cells.Item(131,14).Value = 4.33
$ws.Cells.Item(131,15).Value = "10/11/2023 02:13"
$ws.Cells.Item(131,16).Value = 4.28
$ws.Cells.Item(131,17).Value = "11/11/2023 13:57"
$ws.Cells.Item(131,18).Value = 5.51
$ws.Cells.Item(131,19).Value = "10/11/2023 02:13"
$ws.Cells.Item(131,20).Value = 4.69
$ws.Cells.Item(131,21).Value = "11/11/2023 13:57"
$ws.Cells.Item(131,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-fuzinar/zsHEmc9K/"

# new row 132
$ws.Cells.Item(129,1).Copy()
$ws.Cells.Item(132,1).PasteSpecial(-4122)
$ws.Cells.Item(129,5).Copy()
$ws.Cells.Item(132,5).PasteSpecial(-4122)
$ws.Cells.Item(132,1).Value = 131
$ws.Cells.Item(132,2).Value = "slovenia"
$ws.Cells.Item(132,3).Value = "2-snl"
$ws.Cells.Item(132,4).Value = "2023-2024"
$ws.Cells.Item(132,5).Value = 45241.58333333334
$ws.Cells.Item(132,6).Value = "Bilje"
$ws.Cells.Item(132,7).Value = 0
$ws.Cells.Item(132,8).Value = "ND Gorica"
$ws.Cells.Item(132,9).Value = 1
$ws.Cells.Item(132,10).Value = 3.31
$ws.Cells.Item(132,11).Value = "10/11/2023 02:13"
$ws.Cells.Item(132,12).Value = 4.22
$ws.Cells.Item(132,13).Value = "11/11/2023 13:59"
$ws.Cells.Item(132,14).Value = 3.41
$ws.Cells.Item(132,15).Value = "10/11/2023 02:13"
$ws.Cells.Item(132,16).Value = 4.12
$ws.Cells.Item(132,17).Value = "11/11/2023 13:59"
$ws.Cells.Item(132,18).Value = 1.87
$ws.Cells.Item(132,19).Value = "10/11/2023 02:13"
$ws.Cells.Item(132,20).Value = 1.65
$ws.Cells.Item(132,21).Value = "11/11/2023 13:59"
$ws.Cells.Item(132,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-nd-gorica/tKHAlwfE/"

# new row 133
$ws.Cells.Item(129,1).Copy()
$ws.Cells.Item(133,1).PasteSpecial(-4122)
$ws.Cells.Item(129,5).Copy()
$ws.Cells.Item(133,5).PasteSpecial(-4122)
$ws.Cells.Item(133,1).Value = 132
$ws.Cells.Item(133,2).Value = "slovenia"
$ws.Cells.Item(133,3).Value = "2-snl"
$ws.Cells.Item(133,4).Value = "2023-2024"
$ws.Cells.Item(133,5).Value = 45241.58333333334
$ws.Cells.Item(133,6).Value = "Jadran Dekani"
$ws.Cells.Item(133,7).Value = 3
$ws.Cells.Item(133,8).Value = "NK Krka"
$ws.Cells.Item(133,9).Value = 2
$ws.Cells.Item(133,10).Value = 2.35
$ws.Cells.Item(133,11).Value = "10/11/2023 02:13"
$ws.Cells.Item(133,12).Value = 2.41
$ws.Cells.Item(133,13).Value = "11/11/2023 12:17"
$ws.Cells.Item(133,14).Value = 3.16
$ws.Cells.Item(133,15).Value = "10/11/2023 02:13"
$ws.Cells.Item(133,16).Value = 3.46
$ws.Cells.Item(133,17).Value = "11/11/2023 12:17"
$ws.Cells.Item(133,18).Value = 2.6
$ws.Cells.Item(133,19).Value = "10/11/2023 02:13"
$ws.Cells.Item(133,20).Value = 2.63
$ws.Cells.Item(133,21).Value = "11/11/2023 12:17"
$ws.Cells.Item(133,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/jadran-dekani-nk-krka/jkwAkyP1/"

$excel.CutCopyMode = $false

